$wb = $excel.ActiveWorkbook

function Add-SheetAtEnd {
    param($name)
    $count = $wb.Worksheets.Count
    $lastSheet = $wb.Worksheets.Item($count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $name
    return $newSheet
}

# Add the three new worksheets, in order, after "Message"
$cypherOutputMessage = Add-SheetAtEnd 'CypherOutput_Message'
$statOutput = Add-SheetAtEnd 'StatOutput'
$statOutputMessage = Add-SheetAtEnd 'StatOutput_Message'

# --- CypherOutput_Message: identical layout/content to the "Message" sheet ---
$cypherOutputMessage.Range("A1").Value = 'Neo4j_URL:'
$cypherOutputMessage.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$cypherOutputMessage.Range("A3").Value = 'User_name:'
$cypherOutputMessage.Range("A4").Value = 'neo4j'
$cypherOutputMessage.Range("A5").Value = 'PWD:'
$cypherOutputMessage.Range("A6").Value = 'icdcDBneo4j0'
$cypherOutputMessage.Range("A7").Value = 'Cypher:'
$cypherOutputMessage.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Stage 3''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$cypherOutputMessage.Range("A9").Value = 'Output:'
$cypherOutputMessage.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC18_Canine_Filter_Diagnosis-RespCarciStg3_Neo4jData.xlsx'

# --- StatOutput: header row + one data row of stats ---
$statOutput.Range("A1").Value = 'number_of_files'
$statOutput.Range("B1").Value = 'number_of_sample'
$statOutput.Range("C1").Value = 'number_of_cases'
$statOutput.Range("D1").Value = 'number_of_study'
# (leading apostrophe forces these numeric-looking values to be stored as text, like the source)
$statOutput.Range("A2").Value = '''0'
$statOutput.Range("B2").Value = '''0'
$statOutput.Range("C2").Value = '''38'
$statOutput.Range("D2").Value = '''1'

# --- StatOutput_Message: Message layout repeated twice, second copy uses the StatOutput cypher query ---
$statOutputMessage.Range("A1").Value = 'Neo4j_URL:'
$statOutputMessage.Range("A2").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$statOutputMessage.Range("A3").Value = 'User_name:'
$statOutputMessage.Range("A4").Value = 'neo4j'
$statOutputMessage.Range("A5").Value = 'PWD:'
$statOutputMessage.Range("A6").Value = 'icdcDBneo4j0'
$statOutputMessage.Range("A7").Value = 'Cypher:'
$statOutputMessage.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN [''Stage 3''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$statOutputMessage.Range("A9").Value = 'Output:'
$statOutputMessage.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC18_Canine_Filter_Diagnosis-RespCarciStg3_Neo4jData.xlsx'
$statOutputMessage.Range("A11").Value = 'Neo4j_URL:'
$statOutputMessage.Range("A12").Value = 'bolt://ncias-q2251-c.nci.nih.gov:7687'
$statOutputMessage.Range("A13").Value = 'User_name:'
$statOutputMessage.Range("A14").Value = 'neo4j'
$statOutputMessage.Range("A15").Value = 'PWD:'
$statOutputMessage.Range("A16").Value = 'icdcDBneo4j0'
$statOutputMessage.Range("A17").Value = 'Cypher:'
$statOutputMessage.Range("A18").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN[''Stage 3'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$statOutputMessage.Range("A19").Value = 'Output:'
$statOutputMessage.Range("A20").Value = 'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC18_Canine_Filter_Diagnosis-RespCarciStg3_Neo4jData.xlsx'

# Keep the first sheet ("CypherOutput") as the selected/active tab, matching the source
$wb.Worksheets.Item(1).Activate()

